$p = $ppt.ActivePresentation

# --- Slide 6: Title "Results - LM1" -> "Model Results - LM1" ---
$s6 = $p.Slides.Item(6)
$s6.Shapes.Item(1).TextFrame.TextRange.Text = "Model Results - LM1"

# --- Slide 8: move the architecture diagram picture down slightly ---
$s8 = $p.Slides.Item(8)
$pic = $s8.Shapes.Item(3)
$pic.Top = (1388127 / 12700.0) + 0.00005

# --- Slide 9: Conclusions bullet text tweaks ---
$s9 = $p.Slides.Item(9)
$tr = $s9.Shapes.Item(2).TextFrame.TextRange

# Paragraph 4: "... have to be factored." -> "... post June 2028 should be factored."
$para4 = $tr.Paragraphs(4, 1)
$run4 = $para4.Runs(1, 1)
$run4.Text = "Growth rate assumptions for macroeconomic indicators post June 2028 should be factored."

# Paragraph 6: add a comma after "scale"
$para6 = $tr.Paragraphs(6, 1)
$run6 = $para6.Runs(1, 1)
$run6.Text = "Proposed architecture is easy to scale, and medallion architecture is industry standard."

# Paragraph 7: "Mlflow" -> "MLFlow" and drop the stray trailing quote
$para7 = $tr.Paragraphs(7, 1)
$mlflow = $tr.Characters($para7.Start + 7, 6)
$mlflow.Text = "MLFlow"

$para7b = $tr.Paragraphs(7, 1)
$trailingQuote = $tr.Characters($para7b.Start + $para7b.Length - 2, 1)
$trailingQuote.Text = ""
